$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the "Absent" column (H) values to complete the consolidated report.
$ws.Range("H3").Value = 1
$ws.Range("H6").Value = 0
$ws.Range("H9").Value = 1
$ws.Range("H12").Value = 0
$ws.Range("H15").Value = 1
$ws.Range("H16").Value = 0
